$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "2 5 5"
$ws.Range("B1").Value = "3 7 7"
$ws.Range("C1").ClearContents()

$ws.Range("B2").Select()
